# Swap the contents of columns B:AD between specific row pairs.
# Column A (the running id) stays put; all other columns (B..AD) are
# exchanged between the two rows of each pair, as described by the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pairs = @(
    @(26, 27),
    @(43, 44),
    @(61, 62),
    @(73, 74),
    @(173, 174),
    @(177, 178),
    @(194, 195),
    @(201, 202)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    # Columns B (2) through AD (30)
    for ($col = 2; $col -le 30; $col++) {
        $cell1 = $ws.Cells.Item($r1, $col)
        $cell2 = $ws.Cells.Item($r2, $col)

        $temp = $cell1.Value2
        $cell1.Value2 = $cell2.Value2
        $cell2.Value2 = $temp
    }
}
